$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Electrode Locations" header in C1, matching the existing header style
$ws.Cells.Item(1, 3).Value = "Electrode Locations"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

$data = @(
    @("A1_monopolar_10V_100kHz.txt", 16.345631, "A1"),
    @("A4_monopolar_10V_100kHz.txt", 15.118662, "A4"),
    @("A5_monopolar_10V_100kHz.txt", 13.935383, "A5"),
    @("A7_monopolar_10V_100kHz.txt", 14.181141, "A7"),
    @("A9_monopolar_10V_100kHz.txt", 12.77577, "A9"),
    @("A11_monopolar_10V_100kHz.txt", 9.54997, "A11"),
    @("A13_monopolar_10V_100kHz.txt", 6.014697, "A13"),
    @("B3_monopolar_10V_100kHz.txt", 15.522797, "B3"),
    @("B15_monopolar_10V_100kHz.txt", 9.167680000000001, "B15"),
    @("C1_monopolar_10V_100kHz.txt", 17.385096, "C1"),
    @("C3_monopolar_10V_100kHz.txt", 15.883242, "C3"),
    @("C5_monopolar_10V_100kHz.txt", 13.949946, "C5"),
    @("C7_monopolar_10V_100kHz.txt", 13.762442, "C7"),
    @("C9_monopolar_10V_100kHz.txt", 13.110729, "C9"),
    @("C11_monopolar_10V_100kHz.txt", 10.647688, "C11"),
    @("C14_monopolar_10V_100kHz.txt", 10.431057, "C14"),
    @("C15_monopolar_10V_100kHz.txt", 8.654318999999999, "C15"),
    @("E2_monopolar_10V_100kHz.txt", 17.842023, "E2"),
    @("E3_monopolar_10V_100kHz.txt", 15.806784, "E3"),
    @("E5_monopolar_10V_100kHz.txt", 15.018538, "E5"),
    @("E7_monopolar_10V_100kHz.txt", 14.978489, "E7"),
    @("E9_monopolar_10V_100kHz.txt", 13.656857, "E9"),
    @("E11_monopolar_10V_100kHz.txt", 13.62591, "E11"),
    @("E13_monopolar_10V_100kHz.txt", 10.855217, "E13"),
    @("E15_monopolar_10V_100kHz.txt", 8.860027000000001, "E15"),
    @("F12_monopolar_10V_100kHz.txt", 11.284838, "F12"),
    @("G1_monopolar_10V_100kHz.txt", 18.091422, "G1"),
    @("G3_monopolar_10V_100kHz.txt", 15.426314, "G3"),
    @("G5_monopolar_10V_100kHz.txt", 15.661149, "G5"),
    @("G7_monopolar_10V_100kHz.txt", 15.138686, "G7"),
    @("G9_monopolar_10V_100kHz.txt", 14.199345, "G9"),
    @("G11_monopolar_10V_100kHz.txt", 13.87895, "G11"),
    @("G13_monopolar_10V_100kHz.txt", 11.144665, "G13"),
    @("G15_monopolar_10V_100kHz.txt", 9.0002, "G15"),
    @("H14_monopolar_10V_100kHz.txt", 11.499649, "H14"),
    @("I1_monopolar_10V_100kHz.txt", 18.386331, "I1"),
    @("I3_monopolar_10V_100kHz.txt", 16.624156, "I3"),
    @("I5_monopolar_10V_100kHz.txt", 15.371701, "I5"),
    @("I6_monopolar_10V_100kHz.txt", 13.467533, "I6"),
    @("I7_monopolar_10V_100kHz.txt", 14.740012, "I7"),
    @("I9_monopolar_10V_100kHz.txt", 14.153834, "I9"),
    @("I11_monopolar_10V_100kHz.txt", 14.332236, "I11"),
    @("I13_monopolar_10V_100kHz.txt", 11.425011, "I13"),
    @("I15_monopolar_10V_100kHz.txt", 9.504459000000001, "I15"),
    @("K1_monopolar_10V_100kHz.txt", 19.247394, "K1"),
    @("K3_monopolar_10V_100kHz.txt", 9.619146000000001, "K3"),
    @("K4_monopolar_10V_100kHz.txt", 12.426247, "K4"),
    @("K5_monopolar_10V_100kHz.txt", 14.492434, "K5"),
    @("K7_monopolar_10V_100kHz.txt", 14.328595, "K7"),
    @("K9_monopolar_10V_100kHz.txt", 14.592558, "K9"),
    @("K12_monopolar_10V_100kHz.txt", 12.684748, "K12"),
    @("K13_monopolar_10V_100kHz.txt", 11.567005, "K13"),
    @("K15_monopolar_10V_100kHz.txt", 9.331518000000001, "K15"),
    @("M1_monopolar_10V_100kHz.txt", 20.286859, "M1"),
    @("M3_monopolar_10V_100kHz.txt", 18.373588, "M3"),
    @("M5_monopolar_10V_100kHz.txt", 16.163588, "M5"),
    @("M7_monopolar_10V_100kHz.txt", 14.774601, "M7"),
    @("M9_monopolar_10V_100kHz.txt", 13.775185, "M9"),
    @("M12_monopolar_10V_100kHz.txt", 11.952936, "M12"),
    @("M14_monopolar_10V_100kHz.txt", 10.585794, "M14"),
    @("M15_monopolar_10V_100kHz.txt", 10.964443, "M15"),
    @("O1_monopolar_10V_100kHz.txt", 16.564082, "O1"),
    @("O3_monopolar_10V_100kHz.txt", 17.412402, "O3"),
    @("O5_monopolar_10V_100kHz.txt", 16.17451, "O5"),
    @("O7_monopolar_10V_100kHz.txt", 15.504593, "O7"),
    @("O9_monopolar_10V_100kHz.txt", 14.47423, "O9"),
    @("O11_monopolar_10V_100kHz.txt", 14.170218, "O11"),
    @("O13_monopolar_10V_100kHz.txt", 12.477219, "O13"),
    @("O15_monopolar_10V_100kHz.txt", 11.992985, "O15")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
